$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q3" (same
#    layout/style as the other quarterly sheets) and placing it right
#    before the "总计" (totals) sheet.
# ------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2021-Q3")
$srcIndex = $srcQ3.Index
$srcQ3.Copy($null, $srcQ3)
$newSheet = $wb.Worksheets.Item($srcIndex + 1)
$newSheet.Name = "2022-Q1"

# Header row - only the "基金规模" (fund size) label differs from the
# other quarterly sheets, which use "基金金额" instead.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Make sure the numeric-looking figures (fund code + the money/ratio
# columns) are kept as text, matching the rest of the workbook.
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

# Row 2 - fund 002379
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002379"
$newSheet.Range("C2").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$newSheet.Range("D2").Value = "1.84"
$newSheet.Range("E2").Value = "86.48"
$newSheet.Range("F2").Value = "3.86"
$newSheet.Range("G2").Value = "0.0710"
$newSheet.Range("H2").Value = 6

# Row 3 - fund 002380
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "002380"
$newSheet.Range("C3").Value = "工银瑞信香港中小盘股票（QDII）美元"
$newSheet.Range("D3").Value = "1.84"
$newSheet.Range("E3").Value = "86.48"
$newSheet.Range("F3").Value = "3.86"
$newSheet.Range("G3").Value = "0.0710"
$newSheet.Range("H3").Value = 6

# ------------------------------------------------------------------
# 2. Add a 2022-Q1 row at the top of the "总计" summary sheet, pushing
#    the existing rows down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing data rows (2-5) down to (3-6). Column A just holds a
# sequential 0-based row index, so its value is recomputed rather than
# copied; its styling (bold/bordered, style index 2) is carried over
# via a format-only paste. Columns B-D carry their values down as-is.
for ($r = 5; $r -ge 2; $r--) {
    $destRow = $r + 1
    $total.Range("A$r").Copy()
    $total.Range("A$destRow").PasteSpecial(-4122)
    $total.Range("A$destRow").Value = $destRow - 2
    $total.Range("B$r`:D$r").Copy()
    $total.Range("B$destRow`:D$destRow").PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# New first data row: 2022-Q1
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.14

# Restore the originally active/selected sheet (creating/copying sheets
# along the way shifts Excel's active-tab selection).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "Final sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name
}
